# Updates the "cryptos" price/volume/hour snapshot table (rows 2-51) to the
# refreshed values from the next GitHub Actions scrape run. Price (D) and
# 1h-volume-change (E) are refreshed for the rows that carry real data, and
# the "Hora" column (G) advances from 10 -> 11 for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; D="314.62"; E="2.95%"; G="11"},
    @{Row=3; D="35.31"; E="-2.48%"; G="11"},
    @{Row=4; D="5.130"; E="0.55%"; G="11"},
    @{Row=5; D="0.08102"; E="2.94%"; G="11"},
    @{Row=6; D="2.127"; E="-0.27%"; G="11"},
    @{Row=7; D="7.997"; E="0.71%"; G="11"},
    @{Row=8; D="4.153"; E="1.04%"; G="11"},
    @{Row=9; D="0.9277"; E="0.83%"; G="11"},
    @{Row=10; D="0.1011"; E="4.68%"; G="11"},
    @{Row=11; D="0.1866"; E="0.39%"; G="11"},
    @{Row=12; D="0.09143"; E="5.35%"; G="11"},
    @{Row=13; D="0.03602"; E="1.25%"; G="11"},
    @{Row=14; D="0.09906"; E="-0.40%"; G="11"},
    @{Row=15; D="0.001431"; E="-0.53%"; G="11"},
    @{Row=16; D="0.005787"; E="2.78%"; G="11"},
    @{Row=17; D="3.461"; E="-0.09%"; G="11"},
    @{Row=18; D="2.878"; E="3.75%"; G="11"},
    @{Row=19; D=$null; E="0.78%"; G="11"},
    @{Row=20; D=$null; E="0.27%"; G="11"},
    @{Row=21; D="5.153"; E="-0.23%"; G="11"},
    @{Row=22; D=$null; E="9.96%"; G="11"},
    @{Row=23; D="0.04564"; E="0.20%"; G="11"},
    @{Row=24; D=$null; E="0.85%"; G="11"},
    @{Row=25; D="0.004702"; E="-6.83%"; G="11"},
    @{Row=26; D="0.0001253"; E="-21.86%"; G="11"},
    @{Row=27; D="0.0004507"; E="-5.05%"; G="11"},
    @{Row=28; D=$null; E=$null; G="11"},
    @{Row=29; D=$null; E=$null; G="11"},
    @{Row=30; D=$null; E=$null; G="11"},
    @{Row=31; D=$null; E=$null; G="11"},
    @{Row=32; D=$null; E=$null; G="11"},
    @{Row=33; D=$null; E=$null; G="11"},
    @{Row=34; D=$null; E=$null; G="11"},
    @{Row=35; D=$null; E=$null; G="11"},
    @{Row=36; D=$null; E=$null; G="11"},
    @{Row=37; D=$null; E=$null; G="11"},
    @{Row=38; D=$null; E=$null; G="11"},
    @{Row=39; D="0.01961"; E="6.12%"; G="11"},
    @{Row=40; D="0.04845"; E="1.77%"; G="11"},
    @{Row=41; D="0.007740"; E="3.33%"; G="11"},
    @{Row=42; D=$null; E="-0.57%"; G="11"},
    @{Row=43; D="0.007839"; E="1.29%"; G="11"},
    @{Row=44; D="0.002144"; E="-4.27%"; G="11"},
    @{Row=45; D=$null; E="2.70%"; G="11"},
    @{Row=46; D="0.00006547"; E="3.38%"; G="11"},
    @{Row=47; D=$null; E="0.21%"; G="11"},
    @{Row=48; D="39.24"; E="-17.32%"; G="11"},
    @{Row=49; D="0.001703"; E="-14.81%"; G="11"},
    @{Row=50; D=$null; E="0.21%"; G="11"},
    @{Row=51; D=$null; E="0.21%"; G="11"}
)

# All source cells are stored as text (price/percent/hour strings), so each
# write briefly forces a Text number format before assigning the value -
# otherwise Excel would auto-coerce strings like "314.62" or "11" into
# numbers - then restores the "Normal" style so formatting is unchanged.
foreach ($item in $data) {
    $r = $item.Row
    if ($item.D -ne $null) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($item.E -ne $null) {
        $cell = $ws.Cells.Item($r, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $item.E
        $cell.Style = "Normal"
    }
    if ($item.G -ne $null) {
        $cell = $ws.Cells.Item($r, 7)
        $cell.NumberFormat = "@"
        $cell.Value = $item.G
        $cell.Style = "Normal"
    }
}
